$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(15, 8).Value = 3797.606  # H15: 3975 -> 3797.606
$ws.Cells.Item(15, 9).Value = 3797.606  # I15: 3975 -> 3797.606
$ws.Cells.Item(15, 11).Value = 11392.818  # K15: 11925 -> 11392.818
$ws.Cells.Item(15, 13).Value = -11223.818  # M15: -11756 -> -11223.818
$ws.Cells.Item(28, 8).Value = 289.85715  # H28: 283.86667 -> 289.85715
$ws.Cells.Item(28, 9).Value = 225.66667  # I28: 223.1 -> 225.66667
$ws.Cells.Item(28, 11).Value = 225.66667  # K28: 223.1 -> 225.66667
$ws.Cells.Item(28, 13).Value = 259.33333  # M28: 261.9 -> 259.33333
$ws.Cells.Item(70, 8).Value = 1480.8334  # H70: 1617 -> 1480.8334
$ws.Cells.Item(70, 9).Value = 1150  # I70: 1500 -> 1150
$ws.Cells.Item(70, 11).Value = 3450  # K70: 4500 -> 3450
$ws.Cells.Item(70, 13).Value = -3180  # M70: -4230 -> -3180
$ws.Cells.Item(73, 8).Value = 1480.8334  # H73: 1617 -> 1480.8334
$ws.Cells.Item(73, 9).Value = 1150  # I73: 1500 -> 1150
$ws.Cells.Item(73, 11).Value = 3450  # K73: 4500 -> 3450
$ws.Cells.Item(73, 13).Value = -2514  # M73: -3564 -> -2514
$ws.Cells.Item(92, 8).Value = 4104519.8  # H92: 3078509.8 -> 4104519.8
$ws.Cells.Item(92, 9).Value = 6156281  # I92: 4104347.2 -> 6156281
$ws.Cells.Item(92, 11).Value = 6156281  # K92: 4104347.2 -> 6156281
$ws.Cells.Item(92, 13).Value = -6155033  # M92: -4103099.2 -> -6155033
$ws.Cells.Item(106, 8).Value = 2813.6667  # H106: 1338.4286 -> 2813.6667
$ws.Cells.Item(106, 9).Value = 2813.6667  # I106: 1338.4286 -> 2813.6667
$ws.Cells.Item(106, 11).Value = 2813.6667  # K106: 1338.4286 -> 2813.6667
$ws.Cells.Item(106, 13).Value = -2182.6667  # M106: -707.4286 -> -2182.6667
$ws.Cells.Item(107, 8).Value = 390.7143  # H107: 992.2 -> 390.7143
$ws.Cells.Item(107, 9).Value = 315.5  # I107: 347.33334 -> 315.5
$ws.Cells.Item(107, 10).Value = 842  # J107: 1959.5 -> 842
$ws.Cells.Item(107, 11).Value = 315.5  # K107: 347.33334 -> 315.5
$ws.Cells.Item(107, 12).Value = 842  # L107: 1959.5 -> 842
$ws.Cells.Item(107, 13).Value = 1604.5  # M107: 1572.66666 -> 1604.5
$ws.Cells.Item(107, 14).Value = -4682  # N107: -5799.5 -> -4682
$ws.Cells.Item(113, 8).Value = 9370.714  # H113: 8879.267 -> 9370.714
$ws.Cells.Item(113, 10).Value = 2500  # J113: 2333 -> 2500
$ws.Cells.Item(113, 12).Value = 2500  # L113: 2333 -> 2500
$ws.Cells.Item(113, 14).Value = -9008  # N113: -8841 -> -9008
$ws.Cells.Item(116, 8).Value = 6000  # H116: 16443.111 -> 6000
$ws.Cells.Item(116, 9).Value = 0  # I116: 100000 -> 0
$ws.Cells.Item(116, 10).Value = 6000  # J116: 5998.5 -> 6000
$ws.Cells.Item(116, 11).Value = 0  # K116: 100000 -> 0
$ws.Cells.Item(116, 12).Value = 6000  # L116: 5998.5 -> 6000
$ws.Cells.Item(116, 13).Value = ""  # M116: was -96558 -> removed
$ws.Cells.Item(116, 14).Value = -12884  # N116: -12882.5 -> -12884
$ws.Cells.Item(125, 8).Value = 1135.8334  # H125: 975.5 -> 1135.8334
$ws.Cells.Item(125, 9).Value = 1650  # I125: 1266.6666 -> 1650
$ws.Cells.Item(125, 10).Value = 878.75  # J125: 800.8 -> 878.75
$ws.Cells.Item(125, 11).Value = 14850  # K125: 11399.9994 -> 14850
$ws.Cells.Item(125, 12).Value = 7908.75  # L125: 7207.2 -> 7908.75
$ws.Cells.Item(125, 13).Value = -12390  # M125: -8939.999400000001 -> -12390
$ws.Cells.Item(125, 14).Value = -12828.75  # N125: -12127.2 -> -12828.75
$ws.Cells.Item(132, 8).Value = 990.23254  # H132: 969.9773 -> 990.23254
$ws.Cells.Item(132, 9).Value = 897.0732  # I132: 878.0714 -> 897.0732
$ws.Cells.Item(132, 11).Value = 2691.2196  # K132: 2634.2142 -> 2691.2196
$ws.Cells.Item(132, 13).Value = -161.2196000000004  # M132: -104.2142000000003 -> -161.2196000000004
$ws.Cells.Item(138, 8).Value = 3526.5806  # H138: 3227.3157 -> 3526.5806
$ws.Cells.Item(138, 10).Value = 2245.48  # J138: 2170.3438 -> 2245.48
$ws.Cells.Item(138, 12).Value = 6736.440000000001  # L138: 6511.0314 -> 6736.440000000001
$ws.Cells.Item(138, 14).Value = -17016.44  # N138: -16791.0314 -> -17016.44
$ws.Cells.Item(140, 8).Value = 49023.668  # H140: 81835.31 -> 49023.668
$ws.Cells.Item(140, 10).Value = 49023.668  # J140: 81835.31 -> 49023.668
$ws.Cells.Item(140, 12).Value = 49023.668  # L140: 81835.31 -> 49023.668
$ws.Cells.Item(140, 14).Value = -59383.668  # N140: -92195.31 -> -59383.668
$ws.Cells.Item(141, 8).Value = 905682.6  # H141: 2871.276 -> 905682.6
$ws.Cells.Item(141, 9).Value = 1078921.6  # I141: 2210.72 -> 1078921.6
$ws.Cells.Item(141, 10).Value = 4840  # J141: 6999.75 -> 4840
$ws.Cells.Item(141, 11).Value = 3236764.8  # K141: 6632.16 -> 3236764.8
$ws.Cells.Item(141, 12).Value = 14520  # L141: 20999.25 -> 14520
$ws.Cells.Item(141, 13).Value = -3231584.8  # M141: -1452.16 -> -3231584.8
$ws.Cells.Item(141, 14).Value = -24880  # N141: -31359.25 -> -24880

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 2561.4285  # H32: 2388.2727 -> 2561.4285
$ws.Cells.Item(32, 9).Value = 1808.2878  # I32: 1582.7632 -> 1808.2878
$ws.Cells.Item(32, 10).Value = 7080.273  # J32: 7489.8335 -> 7080.273
$ws.Cells.Item(32, 11).Value = 1808.2878  # K32: 1582.7632 -> 1808.2878
$ws.Cells.Item(32, 12).Value = 7080.273  # L32: 7489.8335 -> 7080.273
$ws.Cells.Item(32, 13).Value = -1521.2878  # M32: -1295.7632 -> -1521.2878
$ws.Cells.Item(32, 14).Value = -7654.273  # N32: -8063.8335 -> -7654.273
$ws.Cells.Item(74, 8).Value = 1579.5217  # H74: 1596.2273 -> 1579.5217
$ws.Cells.Item(74, 9).Value = 1443.5  # I74: 1467.3334 -> 1443.5
$ws.Cells.Item(74, 10).Value = 1727.909  # J74: 1750.9 -> 1727.909
$ws.Cells.Item(74, 11).Value = 1443.5  # K74: 1467.3334 -> 1443.5
$ws.Cells.Item(74, 12).Value = 1727.909  # L74: 1750.9 -> 1727.909
$ws.Cells.Item(74, 13).Value = -569.5  # M74: -593.3334 -> -569.5
$ws.Cells.Item(74, 14).Value = -3475.909  # N74: -3498.9 -> -3475.909
$ws.Cells.Item(77, 8).Value = 1579.5217  # H77: 1596.2273 -> 1579.5217
$ws.Cells.Item(77, 9).Value = 1443.5  # I77: 1467.3334 -> 1443.5
$ws.Cells.Item(77, 10).Value = 1727.909  # J77: 1750.9 -> 1727.909
$ws.Cells.Item(77, 11).Value = 7217.5  # K77: 7336.666999999999 -> 7217.5
$ws.Cells.Item(77, 12).Value = 8639.545  # L77: 8754.5 -> 8639.545
$ws.Cells.Item(77, 13).Value = -2849.5  # M77: -2968.666999999999 -> -2849.5
$ws.Cells.Item(77, 14).Value = -17375.545  # N77: -17490.5 -> -17375.545
$ws.Cells.Item(102, 8).Value = 2398.0908  # H102: 2563.625 -> 2398.0908
$ws.Cells.Item(102, 9).Value = 2153.2222  # I102: 2251.5 -> 2153.2222
$ws.Cells.Item(102, 11).Value = 2153.2222  # K102: 2251.5 -> 2153.2222
$ws.Cells.Item(102, 13).Value = -531.2222000000002  # M102: -629.5 -> -531.2222000000002
$ws.Cells.Item(109, 8).Value = 46444  # H109: 44996 -> 46444
$ws.Cells.Item(109, 10).Value = 46444  # J109: 44996 -> 46444
$ws.Cells.Item(109, 12).Value = 46444  # L109: 44996 -> 46444
$ws.Cells.Item(109, 14).Value = -49218  # N109: -47770 -> -49218
$ws.Cells.Item(110, 8).Value = 4304.3335  # H110: 3553.25 -> 4304.3335
$ws.Cells.Item(110, 9).Value = 1450  # I110: 1400 -> 1450
$ws.Cells.Item(110, 11).Value = 1450  # K110: 1400 -> 1450
$ws.Cells.Item(110, 13).Value = 595  # M110: 645 -> 595
$ws.Cells.Item(132, 8).Value = 1588.7805  # H132: 1415.804 -> 1588.7805
$ws.Cells.Item(132, 9).Value = 1148.4  # I132: 1053.3684 -> 1148.4
$ws.Cells.Item(132, 10).Value = 2789.818  # J132: 2475.2307 -> 2789.818
$ws.Cells.Item(132, 11).Value = 3445.2  # K132: 3160.1052 -> 3445.2
$ws.Cells.Item(132, 12).Value = 8369.454000000002  # L132: 7425.6921 -> 8369.454000000002
$ws.Cells.Item(132, 13).Value = -915.2000000000003  # M132: -630.1052 -> -915.2000000000003
$ws.Cells.Item(132, 14).Value = -13429.454  # N132: -12485.6921 -> -13429.454

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 1460.6666  # H20: 1592.8 -> 1460.6666
$ws.Cells.Item(20, 9).Value = 1582.2142  # I20: 1795.5454 -> 1582.2142
$ws.Cells.Item(20, 11).Value = 1582.2142  # K20: 1795.5454 -> 1582.2142
$ws.Cells.Item(20, 13).Value = -1335.2142  # M20: -1548.5454 -> -1335.2142
$ws.Cells.Item(99, 8).Value = 1765.5555  # H99: 1549.3572 -> 1765.5555
$ws.Cells.Item(99, 9).Value = 1481.6666  # I99: 1243.4445 -> 1481.6666
$ws.Cells.Item(99, 10).Value = 2333.3333  # J99: 2100 -> 2333.3333
$ws.Cells.Item(99, 11).Value = 1481.6666  # K99: 1243.4445 -> 1481.6666
$ws.Cells.Item(99, 12).Value = 2333.3333  # L99: 2100 -> 2333.3333
$ws.Cells.Item(99, 13).Value = 16.33339999999998  # M99: 254.5554999999999 -> 16.33339999999998
$ws.Cells.Item(99, 14).Value = -5329.3333  # N99: -5096 -> -5329.3333
$ws.Cells.Item(107, 8).Value = 4757.75  # H107: 3098.2 -> 4757.75
$ws.Cells.Item(107, 9).Value = 7500.5  # I107: 3562.4285 -> 7500.5
$ws.Cells.Item(107, 11).Value = 7500.5  # K107: 3562.4285 -> 7500.5
$ws.Cells.Item(107, 13).Value = -5580.5  # M107: -1642.4285 -> -5580.5
$ws.Cells.Item(134, 8).Value = 2440.16  # H134: 2014.4688 -> 2440.16
$ws.Cells.Item(134, 9).Value = 2504.7273  # I134: 2019.4138 -> 2504.7273
$ws.Cells.Item(134, 11).Value = 7514.1819  # K134: 6058.2414 -> 7514.1819
$ws.Cells.Item(134, 13).Value = -4979.1819  # M134: -3523.2414 -> -4979.1819

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 1415.8948  # H58: 1550.1471 -> 1415.8948
$ws.Cells.Item(58, 9).Value = 821.7037  # I58: 899.4583 -> 821.7037
$ws.Cells.Item(58, 10).Value = 2874.3635  # J58: 3111.8 -> 2874.3635
$ws.Cells.Item(58, 11).Value = 821.7037  # K58: 899.4583 -> 821.7037
$ws.Cells.Item(58, 12).Value = 2874.3635  # L58: 3111.8 -> 2874.3635
$ws.Cells.Item(58, 13).Value = -618.7037  # M58: -696.4583 -> -618.7037
$ws.Cells.Item(58, 14).Value = -3280.3635  # N58: -3517.8 -> -3280.3635
$ws.Cells.Item(62, 8).Value = 2710.8  # H62: 2261.4443 -> 2710.8
$ws.Cells.Item(62, 9).Value = 2651  # I62: 2701.3333 -> 2651
$ws.Cells.Item(62, 10).Value = 2950  # J62: 2041.5 -> 2950
$ws.Cells.Item(62, 11).Value = 2651  # K62: 2701.3333 -> 2651
$ws.Cells.Item(62, 12).Value = 2950  # L62: 2041.5 -> 2950
$ws.Cells.Item(62, 13).Value = -2027  # M62: -2077.3333 -> -2027
$ws.Cells.Item(62, 14).Value = -4198  # N62: -3289.5 -> -4198
$ws.Cells.Item(65, 8).Value = 2710.8  # H65: 2261.4443 -> 2710.8
$ws.Cells.Item(65, 9).Value = 2651  # I65: 2701.3333 -> 2651
$ws.Cells.Item(65, 10).Value = 2950  # J65: 2041.5 -> 2950
$ws.Cells.Item(65, 11).Value = 13255  # K65: 13506.6665 -> 13255
$ws.Cells.Item(65, 12).Value = 14750  # L65: 10207.5 -> 14750
$ws.Cells.Item(65, 13).Value = -10135  # M65: -10386.6665 -> -10135
$ws.Cells.Item(65, 14).Value = -20990  # N65: -16447.5 -> -20990
$ws.Cells.Item(105, 8).Value = 2302  # H105: 2043.8 -> 2302
$ws.Cells.Item(105, 9).Value = 2252.5  # I105: 1929.75 -> 2252.5
$ws.Cells.Item(105, 11).Value = 2252.5  # K105: 1929.75 -> 2252.5
$ws.Cells.Item(105, 13).Value = -505.5  # M105: -182.75 -> -505.5
$ws.Cells.Item(132, 8).Value = 2281.543  # H132: 2444.7334 -> 2281.543
$ws.Cells.Item(132, 9).Value = 1522.1  # I132: 1560.5883 -> 1522.1
$ws.Cells.Item(132, 10).Value = 3294.1333  # J132: 3600.923 -> 3294.1333
$ws.Cells.Item(132, 11).Value = 4566.299999999999  # K132: 4681.7649 -> 4566.299999999999
$ws.Cells.Item(132, 12).Value = 9882.3999  # L132: 10802.769 -> 9882.3999
$ws.Cells.Item(132, 13).Value = -2036.299999999999  # M132: -2151.7649 -> -2036.299999999999
$ws.Cells.Item(132, 14).Value = -14942.3999  # N132: -15862.769 -> -14942.3999
$ws.Cells.Item(136, 8).Value = 1415.8948  # H136: 1550.1471 -> 1415.8948
$ws.Cells.Item(136, 9).Value = 821.7037  # I136: 899.4583 -> 821.7037
$ws.Cells.Item(136, 10).Value = 2874.3635  # J136: 3111.8 -> 2874.3635
$ws.Cells.Item(136, 11).Value = 2465.1111  # K136: 2698.3749 -> 2465.1111
$ws.Cells.Item(136, 12).Value = 8623.0905  # L136: 9335.400000000001 -> 8623.0905
$ws.Cells.Item(136, 13).Value = 84.88889999999992  # M136: -148.3748999999998 -> 84.88889999999992
$ws.Cells.Item(136, 14).Value = -13723.0905  # N136: -14435.4 -> -13723.0905

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 772.8200000000001  # H131: 793.83 -> 772.8200000000001
$ws.Cells.Item(131, 10).Value = 788.8936  # J131: 811.2447 -> 788.8936
$ws.Cells.Item(131, 12).Value = 2366.6808  # L131: 2433.7341 -> 2366.6808
$ws.Cells.Item(131, 14).Value = -12446.6808  # N131: -12513.7341 -> -12446.6808

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(46, 8).Value = 23499.9  # H46: 27914.143 -> 23499.9
$ws.Cells.Item(46, 10).Value = 23499.9  # J46: 27914.143 -> 23499.9
$ws.Cells.Item(46, 12).Value = 23499.9  # L46: 27914.143 -> 23499.9
$ws.Cells.Item(46, 14).Value = -23811.9  # N46: -28226.143 -> -23811.9
$ws.Cells.Item(70, 8).Value = 10745.363  # H70: 12666 -> 10745.363
$ws.Cells.Item(70, 9).Value = 14528.571  # I70: 19199 -> 14528.571
$ws.Cells.Item(70, 10).Value = 4124.75  # J70: 4499.75 -> 4124.75
$ws.Cells.Item(70, 11).Value = 14528.571  # K70: 19199 -> 14528.571
$ws.Cells.Item(70, 12).Value = 4124.75  # L70: 4499.75 -> 4124.75
$ws.Cells.Item(70, 13).Value = -14258.571  # M70: -18929 -> -14258.571
$ws.Cells.Item(70, 14).Value = -4664.75  # N70: -5039.75 -> -4664.75
$ws.Cells.Item(73, 8).Value = 10745.363  # H73: 12666 -> 10745.363
$ws.Cells.Item(73, 9).Value = 14528.571  # I73: 19199 -> 14528.571
$ws.Cells.Item(73, 10).Value = 4124.75  # J73: 4499.75 -> 4124.75
$ws.Cells.Item(73, 11).Value = 14528.571  # K73: 19199 -> 14528.571
$ws.Cells.Item(73, 12).Value = 4124.75  # L73: 4499.75 -> 4124.75
$ws.Cells.Item(73, 13).Value = -13592.571  # M73: -18263 -> -13592.571
$ws.Cells.Item(73, 14).Value = -5996.75  # N73: -6371.75 -> -5996.75
$ws.Cells.Item(97, 8).Value = 1774.0834  # H97: 2010.1 -> 1774.0834
$ws.Cells.Item(97, 9).Value = 2101.6667  # I97: 2402 -> 2101.6667
$ws.Cells.Item(97, 10).Value = 1446.5  # J97: 1618.2 -> 1446.5
$ws.Cells.Item(97, 11).Value = 2101.6667  # K97: 2402 -> 2101.6667
$ws.Cells.Item(97, 12).Value = 1446.5  # L97: 1618.2 -> 1446.5
$ws.Cells.Item(97, 13).Value = -1605.6667  # M97: -1906 -> -1605.6667
$ws.Cells.Item(97, 14).Value = -2438.5  # N97: -2610.2 -> -2438.5
$ws.Cells.Item(113, 8).Value = 1127  # H113: 1211.5 -> 1127
$ws.Cells.Item(113, 9).Value = 798  # I113: 795 -> 798
$ws.Cells.Item(113, 10).Value = 1236.6666  # J113: 1628 -> 1236.6666
$ws.Cells.Item(113, 11).Value = 798  # K113: 795 -> 798
$ws.Cells.Item(113, 12).Value = 1236.6666  # L113: 1628 -> 1236.6666
$ws.Cells.Item(113, 13).Value = 1372  # M113: 1375 -> 1372
$ws.Cells.Item(113, 14).Value = -5576.6666  # N113: -5968 -> -5576.6666

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(14, 8).Value = 11999  # H14: 10249.75 -> 11999
$ws.Cells.Item(14, 10).Value = 11999  # J14: 10249.75 -> 11999
$ws.Cells.Item(14, 12).Value = 11999  # L14: 10249.75 -> 11999
$ws.Cells.Item(14, 14).Value = -12343  # N14: -10593.75 -> -12343
$ws.Cells.Item(16, 8).Value = 5924.7144  # H16: 5697 -> 5924.7144
$ws.Cells.Item(16, 10).Value = 1894.8  # J16: 2262.8333 -> 1894.8
$ws.Cells.Item(16, 12).Value = 1894.8  # L16: 2262.8333 -> 1894.8
$ws.Cells.Item(16, 14).Value = -2234.8  # N16: -2602.8333 -> -2234.8
$ws.Cells.Item(40, 8).Value = 4608.933  # H40: 4345.1875 -> 4608.933
$ws.Cells.Item(40, 9).Value = 2053  # I40: 1868.1111 -> 2053
$ws.Cells.Item(40, 11).Value = 2053  # K40: 1868.1111 -> 2053
$ws.Cells.Item(40, 13).Value = -1917  # M40: -1732.1111 -> -1917
$ws.Cells.Item(132, 8).Value = 2366.1943  # H132: 2447.0293 -> 2366.1943
$ws.Cells.Item(132, 9).Value = 1938.1875  # I132: 2073.3572 -> 1938.1875
$ws.Cells.Item(132, 11).Value = 5814.5625  # K132: 6220.071599999999 -> 5814.5625
$ws.Cells.Item(132, 13).Value = -3284.5625  # M132: -3690.071599999999 -> -3284.5625
$ws.Cells.Item(136, 8).Value = 3157.5  # H136: 3178.2173 -> 3157.5
$ws.Cells.Item(136, 9).Value = 2605.1  # I136: 2678.818 -> 2605.1
$ws.Cells.Item(136, 10).Value = 3617.8333  # J136: 3636 -> 3617.8333
$ws.Cells.Item(136, 11).Value = 7815.299999999999  # K136: 8036.454000000001 -> 7815.299999999999
$ws.Cells.Item(136, 12).Value = 10853.4999  # L136: 10908 -> 10853.4999
$ws.Cells.Item(136, 13).Value = -5265.299999999999  # M136: -5486.454000000001 -> -5265.299999999999
$ws.Cells.Item(136, 14).Value = -15953.4999  # N136: -16008 -> -15953.4999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 56677  # H122: 34754.434 -> 56677
$ws.Cells.Item(122, 9).Value = 60929.152  # I122: 37979.19 -> 60929.152
$ws.Cells.Item(122, 10).Value = 1399  # J122: 894.5 -> 1399
$ws.Cells.Item(122, 11).Value = 182787.456  # K122: 113937.57 -> 182787.456
$ws.Cells.Item(122, 12).Value = 4197  # L122: 2683.5 -> 4197
$ws.Cells.Item(122, 13).Value = -180337.456  # M122: -111487.57 -> -180337.456
$ws.Cells.Item(122, 14).Value = -9097  # N122: -7583.5 -> -9097
$ws.Cells.Item(123, 8).Value = 45306.062  # H123: 45166.668 -> 45306.062
$ws.Cells.Item(123, 10).Value = 47499.785  # J123: 47507.69 -> 47499.785
$ws.Cells.Item(123, 12).Value = 47499.785  # L123: 47507.69 -> 47499.785
$ws.Cells.Item(123, 14).Value = -57299.785  # N123: -57307.69 -> -57299.785
$ws.Cells.Item(124, 8).Value = 19749  # H124: 0 -> 19749
$ws.Cells.Item(124, 10).Value = 19749  # J124: 0 -> 19749
$ws.Cells.Item(124, 12).Value = 19749  # L124: 0 -> 19749
$ws.Cells.Item(124, 14).Value = -29569  # N124: None -> -29569
$ws.Cells.Item(125, 8).Value = 39978.95  # H125: 40000 -> 39978.95
$ws.Cells.Item(125, 10).Value = 39978.95  # J125: 40000 -> 39978.95
$ws.Cells.Item(125, 12).Value = 39978.95  # L125: 40000 -> 39978.95
$ws.Cells.Item(125, 14).Value = -49818.95  # N125: -49840 -> -49818.95
$ws.Cells.Item(132, 8).Value = 1259.5641  # H132: 1331.3611 -> 1259.5641
$ws.Cells.Item(132, 9).Value = 986.4286  # I132: 1046.9615 -> 986.4286
$ws.Cells.Item(132, 10).Value = 1954.8182  # J132: 2070.8 -> 1954.8182
$ws.Cells.Item(132, 11).Value = 2959.2858  # K132: 3140.8845 -> 2959.2858
$ws.Cells.Item(132, 12).Value = 5864.4546  # L132: 6212.400000000001 -> 5864.4546
$ws.Cells.Item(132, 13).Value = -429.2857999999997  # M132: -610.8844999999997 -> -429.2857999999997
$ws.Cells.Item(132, 14).Value = -10924.4546  # N132: -11272.4 -> -10924.4546
